# Updated cryptos list on Fri Sep 29 20:36:48 UTC 2023 with GitHub Actions
#
# Refreshes the scraped price / 1h-volume figures on "Sheet1" of the
# cryptos workbook, and fixes the ordering of a few rows whose rank
# swapped places (PaxDollar/FraxShare and WEMIXToken/Aave).
#
# Some "Price" values are plain decimal numbers (e.g. "0.0621", "66.07")
# that Excel would otherwise silently coerce to a Number when assigned via
# .Value. Those cells are forced to Text format first so the value is
# written back exactly as scraped (matching the rest of the sheet, which
# stores every price/volume cell as a string).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

function Set-TextValue($addr, $val) {
    $ws.Range($addr).NumberFormat = "@"
    $ws.Range($addr).Value = $val
}

# Row 2 - Bitcoin
$ws.Range("D2").Value = "26.918.79"
$ws.Range("E2").Value = "  -0.71%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "1.667.12"
$ws.Range("E3").Value = "  +0.56%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  -0.02%  "

# Row 5 - BNB
Set-TextValue "D5" "215.53"

# Row 6 - XRP
$ws.Range("E6").Value = "  +4.87%  "

# Row 7 - USDC
$ws.Range("E7").Value = "  -0.04%  "

# Row 8 - Dogecoin
Set-TextValue "D8" "0.0621"
$ws.Range("E8").Value = "  +1.04%  "

# Row 9 - Cardano
$ws.Range("E9").Value = "  -0.45%  "

# Row 10 - Solana
Set-TextValue "D10" "20.25"
$ws.Range("E10").Value = "  +2.55%  "

# Row 11 - TRON
Set-TextValue "D11" "0.0895"
$ws.Range("E11").Value = "  +3.57%  "

# Row 12 - WrappedliquidstakedEther2.0
$ws.Range("D12").Value = "1.901.02"
$ws.Range("E12").Value = "  +0.49%  "

# Row 13 - WrappedEther
$ws.Range("D13").Value = "1.665.94"
$ws.Range("E13").Value = "  +0.24%  "

# Row 14 - Polkadot
$ws.Range("E14").Value = "  -0.09%  "

# Row 15 - Polygon
$ws.Range("E15").Value = "  +0.90%  "

# Row 16 - Litecoin
Set-TextValue "D16" "66.07"
$ws.Range("E16").Value = "  +1.45%  "

# Row 17 - WrappedBTC
$ws.Range("D17").Value = "26.909.51"
$ws.Range("E17").Value = "  -0.70%  "

# Row 18 - BitcoinCash
Set-TextValue "D18" "235.03"
$ws.Range("E18").Value = "  -1.45%  "

# Row 19 - Chainlink
$ws.Range("E19").Value = "  +1.25%  "

# Row 20 - ShibaInu
$ws.Range("E20").Value = "  +0.46%  "

# Row 21 - Dai
$ws.Range("E21").Value = "  +0.06%  "

# Row 22 - Uniswap
Set-TextValue "D22" "4.36"
$ws.Range("E22").Value = "  -1.65%  "

# Row 24 - Avalanche
Set-TextValue "D24" "9.12"
$ws.Range("E24").Value = "  -1.19%  "

# Row 25 - Monero
Set-TextValue "D25" "146.02"
$ws.Range("E25").Value = "  +0.26%  "

# Row 26 - Cosmos
Set-TextValue "D26" "7.12"
$ws.Range("E26").Value = "  -0.27%  "

# Row 27 - Stellar
Set-TextValue "D27" "0.114"
$ws.Range("E27").Value = "  +0.78%  "

# Row 28 - EthereumClassic
Set-TextValue "D28" "15.88"
$ws.Range("E28").Value = "  +0.18%  "

# Row 29 - BinanceUSD
$ws.Range("E29").Value = "  +0.06%  "

# Row 30 - Hedera
$ws.Range("E30").Value = "  -0.25%  "

# Row 31 - PancakeSwap
$ws.Range("E31").Value = "  -0.04%  "

# Row 32 - Filecoin
Set-TextValue "D32" "3.36"
$ws.Range("E32").Value = "  +1.90%  "

# Row 33 - Maker
$ws.Range("D33").Value = "1.457.31"
$ws.Range("E33").Value = "  -4.00%  "

# Row 34 - InternetComputer(DFINITY)
$ws.Range("E34").Value = "  +2.56%  "

# Row 35 - LidoDAOToken
Set-TextValue "D35" "1.64"
$ws.Range("E35").Value = "  +1.96%  "

# Row 36 - HuobiToken
$ws.Range("E36").Value = "  -0.26%  "

# Row 37 - ImmutableX
$ws.Range("E37").Value = "  +0.41%  "

# Row 38 - ARBITRUM
Set-TextValue "D38" "0.904"

# Row 39 - VeChain
$ws.Range("E39").Value = "  +0.14%  "

# Row 40 - now FraxShare (was PaxDollar)
$ws.Range("B40").Value = "FraxShare"
$ws.Range("C40").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
Set-TextValue "D40" "5.71"
$ws.Range("E40").Value = "  -4.05%  "

# Row 41 - now PaxDollar (was FraxShare)
$ws.Range("B41").Value = "PaxDollar"
$ws.Range("C41").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
Set-TextValue "D41" "1.00"
$ws.Range("E41").Value = "  +0.04%  "

# Row 42 - MXToken
$ws.Range("E42").Value = "  +0.46%  "

# Row 43 - now Aave (was WEMIXToken)
$ws.Range("B43").Value = "Aave"
$ws.Range("C43").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
Set-TextValue "D43" "65.94"
$ws.Range("E43").Value = "  -0.44%  "

# Row 44 - now WEMIXToken (was Aave)
$ws.Range("B44").Value = "WEMIXToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
Set-TextValue "D44" "0.973"
$ws.Range("E44").Value = "  +5.81%  "

# Row 45 - RocketPoolETH
$ws.Range("D45").Value = "1.808.71"

# Row 46 - TrustWalletToken
Set-TextValue "D46" "0.782"
$ws.Range("E46").Value = "  +0.49%  "

# Row 47 - Quant
$ws.Range("E47").Value = "  +0.62%  "

# Row 48 - RenderToken
$ws.Range("E48").Value = "  +0.75%  "

# Row 49 - BabyDogeCoin
$ws.Range("E49").Value = "  -1.69%  "

# Row 50 - Algorand
$ws.Range("E50").Value = "  +4.42%  "

# Row 51 - Cronos
$ws.Range("E51").Value = "  -0.02%  "
